# Daily attendance processing - 2026-01-08 15:09:07
# Swap the order of recorder names in column G ("Recorded By") wherever
# the cell currently reads "dnasr281@gmail.com, System" so that it reads
# "System, dnasr281@gmail.com" instead.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
